$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.61829200881043
$ws.Range("C2").Value = 5.637866171132563
$ws.Range("D2").Value = 4.23390688931526
$ws.Range("F2").Value = 55.80349110604007
$ws.Range("G2").Value = 3.789626048486827
$ws.Range("I2").Value = 40.47954401922465
$ws.Range("J2").Value = 10.55481807256826
$ws.Range("K2").Value = 15.83711071410278
$ws.Range("L2").Value = 12.12974397727694

$ws.Range("B3").Value = 17.56981555844717
$ws.Range("C3").Value = 5.505059495140196
$ws.Range("D3").Value = 4.255910923883797
$ws.Range("F3").Value = 55.60503640512128
$ws.Range("G3").Value = 3.792773168910352
$ws.Range("I3").Value = 40.37674773694604
$ws.Range("J3").Value = 10.56822686463771
$ws.Range("K3").Value = 15.80464444538398
$ws.Range("L3").Value = 12.14831809996609

$ws.Range("B4").Value = 17.54518279186487
$ws.Range("C4").Value = 5.423795793643277
$ws.Range("D4").Value = 4.270286594171957
$ws.Range("F4").Value = 55.49026954920508
$ws.Range("G4").Value = 3.794806639824577
$ws.Range("I4").Value = 40.31784229832042
$ws.Range("J4").Value = 10.57716374141505
$ws.Range("K4").Value = 15.78894453242371
$ws.Range("L4").Value = 12.16163812840418

$ws.Range("B5").Value = 17.53644463572526
$ws.Range("C5").Value = 5.39080555970054
$ws.Range("D5").Value = 4.276363031391301
$ws.Range("F5").Value = 55.44530333359429
$ws.Range("G5").Value = 3.795660814763967
$ws.Range("I5").Value = 40.29490478979467
$ws.Range("J5").Value = 10.58098289780157
$ws.Range("K5").Value = 15.78361717024827
$ws.Range("L5").Value = 12.16754798659174

$ws.Range("B6").Value = 17.53507242727811
$ws.Range("C6").Value = 5.385336759067936
$ws.Range("D6").Value = 4.277385220570326
$ws.Range("F6").Value = 55.43794614150465
$ws.Range("G6").Value = 3.795804193794668
$ws.Range("I6").Value = 40.29116068635732
$ws.Range("J6").Value = 10.58162778420092
$ws.Range("K6").Value = 15.78279736517324
$ws.Range("L6").Value = 12.16855842001361

$ws.Range("B7").Value = 17.54505967153036
$ws.Range("C7").Value = 5.423350294616931
$ws.Range("D7").Value = 4.270367658567098
$ws.Range("F7").Value = 55.48965579636407
$ws.Range("G7").Value = 3.7948180560815
$ws.Range("I7").Value = 40.31752862547354
$ws.Range("J7").Value = 10.57721452955696
$ws.Range("K7").Value = 15.78886834474101
$ws.Range("L7").Value = 12.16171587968053

$ws.Range("B8").Value = 17.60051830155825
$ws.Range("C8").Value = 5.592048823151222
$ws.Range("D8").Value = 4.241314746753531
$ws.Range("F8").Value = 55.73360629822542
$ws.Range("G8").Value = 3.790690241135796
$ws.Range("I8").Value = 40.44323060271659
$ws.Range("J8").Value = 10.55929554018078
$ws.Range("K8").Value = 15.82504119754474
$ws.Range("L8").Value = 12.13575091358631

$ws.Range("B9").Value = 17.74950801702469
$ws.Range("C9").Value = 5.922754535166892
$ws.Range("D9").Value = 4.191175105655116
$ws.Range("F9").Value = 56.26723243661394
$ws.Range("G9").Value = 3.783393894347879
$ws.Range("I9").Value = 40.72277997357725
$ws.Range("J9").Value = 10.52972669786905
$ws.Range("K9").Value = 15.92927787224048
$ws.Range("L9").Value = 12.1000221827436

$ws.Range("B10").Value = 17.8827510275332
$ws.Range("C10").Value = 6.16267834986824
$ws.Range("D10").Value = 4.158460375764764
$ws.Range("F10").Value = 56.69155247088958
$ws.Range("G10").Value = 3.778514179784309
$ws.Range("I10").Value = 40.94777493330321
$ws.Range("J10").Value = 10.51137896060045
$ws.Range("K10").Value = 16.02569608231279
$ws.Range("L10").Value = 12.08301617320865

$ws.Range("B11").Value = 17.94834107410593
$ws.Range("C11").Value = 6.270585949722046
$ws.Range("D11").Value = 4.144463831412244
$ws.Range("F11").Value = 56.89125705778457
$ws.Range("G11").Value = 3.776397471920878
$ws.Range("I11").Value = 41.05425855660116
$ws.Range("J11").Value = 10.50376119327922
$ws.Range("K11").Value = 16.07374553839643
$ws.Range("L11").Value = 12.07728222744758

$ws.Range("B12").Value = 17.97387602053933
$ws.Range("C12").Value = 6.311225988222525
$ws.Range("D12").Value = 4.139290348049851
$ws.Range("F12").Value = 56.96780923126389
$ws.Range("G12").Value = 3.775610661645944
$ws.Range("I12").Value = 41.09516277053488
$ws.Range("J12").Value = 10.50098100480113
$ws.Range("K12").Value = 16.09253114499709
$ws.Range("L12").Value = 12.07539829335052

$ws.Range("B13").Value = 17.96834589921605
$ws.Range("C13").Value = 6.302483989619115
$ws.Range("D13").Value = 4.140398925127037
$ws.Range("F13").Value = 56.95128158767792
$ws.Range("G13").Value = 3.775779461042575
$ws.Range("I13").Value = 41.08632768381016
$ws.Range("J13").Value = 10.5015751253703
$ws.Range("K13").Value = 16.08845927409909
$ws.Range("L13").Value = 12.07579125988363

$ws.Range("B14").Value = 17.95042798130073
$ws.Range("C14").Value = 6.27393412839561
$ws.Range("D14").Value = 4.144035670104927
$ws.Range("F14").Value = 56.89753663412017
$ws.Range("G14").Value = 3.776332445614641
$ws.Range("I14").Value = 41.05761220062612
$ws.Range("J14").Value = 10.50353037294877
$ws.Range("K14").Value = 16.07527928245304
$ws.Range("L14").Value = 12.07712147869072

$ws.Range("B15").Value = 17.9395430000506
$ws.Range("C15").Value = 6.256416302390091
$ws.Range("D15").Value = 4.146279763254763
$ws.Range("F15").Value = 56.86473622298423
$ws.Range("G15").Value = 3.776673081855587
$ws.Range("I15").Value = 41.04009842449232
$ws.Range("J15").Value = 10.50474161814432
$ws.Range("K15").Value = 16.06728266336008
$ws.Range("L15").Value = 12.07797368540291

$ws.Range("B16").Value = 17.87856318164278
$ws.Range("C16").Value = 6.155597596013555
$ws.Range("D16").Value = 4.159392847804141
$ws.Range("F16").Value = 56.67863300163071
$ws.Range("G16").Value = 3.778654579004268
$ws.Range("I16").Value = 40.94089803281309
$ws.Range("J16").Value = 10.51189143807962
$ws.Range("K16").Value = 16.02263921669521
$ws.Range("L16").Value = 12.08343114604991

$ws.Range("B17").Value = 17.84241618941792
$ws.Range("C17").Value = 6.093399525702401
$ws.Range("D17").Value = 4.167663655914753
$ws.Range("F17").Value = 56.56615307409829
$ws.Range("G17").Value = 3.779896509161952
$ws.Range("I17").Value = 40.88109114692474
$ws.Range("J17").Value = 10.51646405274361
$ws.Range("K17").Value = 15.99631643371346
$ws.Range("L17").Value = 12.08729157987612

$ws.Range("B18").Value = 17.82209504357141
$ws.Range("C18").Value = 6.057511075083412
$ws.Range("D18").Value = 4.172504195856424
$ws.Range("F18").Value = 56.50208829923847
$ws.Range("G18").Value = 3.780620543703867
$ws.Range("I18").Value = 40.84708158814104
$ws.Range("J18").Value = 10.51916270695908
$ws.Range("K18").Value = 15.9815712676186
$ws.Range("L18").Value = 12.08970048990578

$ws.Range("B19").Value = 17.81529585187949
$ws.Range("C19").Value = 6.045341734780738
$ws.Range("D19").Value = 4.174157461306359
$ws.Range("F19").Value = 56.48050634433056
$ws.Range("G19").Value = 3.780867359496094
$ws.Range("I19").Value = 40.83563386584581
$ws.Range("J19").Value = 10.52008821718186
$ws.Range("K19").Value = 15.97664699856362
$ws.Range("L19").Value = 12.09054849346496

$ws.Range("B20").Value = 17.8462156202734
$ws.Range("C20").Value = 6.100032724183166
$ws.Range("D20").Value = 4.166774588371225
$ws.Range("F20").Value = 56.57806170684986
$ws.Range("G20").Value = 3.779763299331374
$ws.Range("I20").Value = 40.8874174435047
$ws.Range("J20").Value = 10.51597019165454
$ws.Range("K20").Value = 15.99907774134727
$ws.Range("L20").Value = 12.08686112559995

$ws.Range("B21").Value = 17.95567213218533
$ws.Range("C21").Value = 6.282326286003426
$ws.Range("D21").Value = 4.142964036276388
$ws.Range("F21").Value = 56.91329789687344
$ws.Range("G21").Value = 3.776169621243807
$ws.Range("I21").Value = 41.06603096879212
$ws.Range("J21").Value = 10.50295323571831
$ws.Range("K21").Value = 16.079134644263
$ws.Range("L21").Value = 12.07672296617243

$ws.Range("B22").Value = 18.03126405574862
$ws.Range("C22").Value = 6.40015188824323
$ws.Range("D22").Value = 4.128140702313669
$ws.Range("F22").Value = 57.13778953277728
$ws.Range("G22").Value = 3.77390682819642
$ws.Range("I22").Value = 41.1861462034044
$ws.Range("J22").Value = 10.4950548415626
$ws.Range("K22").Value = 16.13489110533792
$ws.Range("L22").Value = 12.07177190090725

$ws.Range("B23").Value = 17.9905544743134
$ws.Range("C23").Value = 6.337400287105972
$ws.Range("D23").Value = 4.135984854477828
$ws.Range("F23").Value = 57.01749142059582
$ws.Range("G23").Value = 3.77510669286789
$ws.Range("I23").Value = 41.12173358414415
$ws.Range("J23").Value = 10.49921474200185
$ws.Range("K23").Value = 16.10482275029611
$ws.Range("L23").Value = 12.07426132702366

$ws.Range("B24").Value = 17.8444964626168
$ws.Range("C24").Value = 6.097034256500456
$ws.Range("D24").Value = 4.167176269315704
$ws.Range("F24").Value = 56.5726759381009
$ws.Range("G24").Value = 3.7798234922398
$ws.Range("I24").Value = 40.88455615814617
$ws.Range("J24").Value = 10.51619324882351
$ws.Range("K24").Value = 15.99782814422323
$ws.Range("L24").Value = 12.08705514369139

$ws.Range("B25").Value = 17.70496633893276
$ws.Range("C25").Value = 5.833618809245854
$ws.Range("D25").Value = 4.204012125793562
$ws.Range("F25").Value = 56.11709720493757
$ws.Range("G25").Value = 3.785282883229418
$ws.Range("I25").Value = 40.64366857541275
$ws.Range("J25").Value = 10.53713151373571
$ws.Range("K25").Value = 15.8975578737501
$ws.Range("L25").Value = 12.10806301218556

